# Weekly update: a new week of price data is prepended to the detail table
# (rows shift down by 2; the two oldest rows fall off the bottom of the
# original range and get re-appended at the end of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of this block (row 221), pushing the
# former rows 221:240 down to 223:240, matching the prior row's formatting.
$ws.Rows("221:222").Insert()

# --- New row 221 ---
$ws.Range("A221").Value = 9
$ws.Range("B221").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C221").Value = "Metropolitana"
$ws.Range("D221").Value = 44461
$ws.Range("E221").Value = 13
$ws.Range("F221").Value = 100114014
$ws.Range("G221").Value = "Betarraga"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 3400
$ws.Range("K221").Value = 90
$ws.Range("L221").Value = 100
$ws.Range("M221").Value = 95
$ws.Range("N221").Value = "$/unidad"
$ws.Range("O221").Value = "Región Metropolitana"
$ws.Range("P221").Value = 95
$ws.Range("Q221").Value = 1
$ws.Range("R221").Value = "Hortaliza"

# --- New row 222 ---
$ws.Range("A222").Value = 9
$ws.Range("B222").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C222").Value = "Metropolitana"
$ws.Range("D222").Value = 44461
$ws.Range("E222").Value = 13
$ws.Range("F222").Value = 100114014
$ws.Range("G222").Value = "Betarraga"
$ws.Range("H222").Value = "Sin especificar"
$ws.Range("I222").Value = "Segunda"
$ws.Range("J222").Value = 1600
$ws.Range("K222").Value = 70
$ws.Range("L222").Value = 80
$ws.Range("M222").Value = 75
$ws.Range("N222").Value = "$/unidad"
$ws.Range("O222").Value = "Región Metropolitana"
$ws.Range("P222").Value = 75
$ws.Range("Q222").Value = 1
$ws.Range("R222").Value = "Hortaliza"
